# Sprint Backlog Week2 update
# 1) Updated time spent on tasks during week 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column B: reassign "TBD" rows (16-20) to "Vitor" ---
$ws.Range("B16").Value = "Vitor"
$ws.Range("B17").Value = "Vitor"
$ws.Range("B18").Value = "Vitor"
$ws.Range("B19").Value = "Vitor"
$ws.Range("B20").Value = "Vitor"

# --- Column E (Week 1 "amount remaining") updates ---
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 3
$ws.Range("E17").Value = 2
$ws.Range("E18").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("E20").Value = 2

# --- Column F (Week 2 "amount remaining") updates ---
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 1

# --- Update the active selection to reflect where editing ended (F20) ---
$ws.Range("F20").Select()

$wb.Save()
